$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Insert a new row before the current row 3 (the "MDA / Campaign / age 15-50" row),
# shifting the later rows down by one. The new row will carry the age band and
# coverage values that used to live on row 2.
$ws.Rows("3:3").Insert()

# --- New row 3: MDA / Treatment / Campaign, age 2-15, coverage 0.8 ---
# (these values previously lived on row 2 - min age 2, max age 15, and the P..AD
#  coverage figures of 0.8 that used to sit alongside them)
$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = $ws.Range("C2").Text
$ws.Range("D3").Value = $ws.Range("D2").Text
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 15
$ws.Range("P3").Value = 0.8
$ws.Range("R3").Value = 0.8
$ws.Range("T3").Value = 0.8
$ws.Range("V3").Value = 0.8
$ws.Range("X3").Value = 0.8
$ws.Range("Z3").Value = 0.8
$ws.Range("AB3").Value = 0.8
$ws.Range("AD3").Value = 0.8

# --- Row 2 is now the MDA / Treatment / Campaign age 5-15 row, MDA-only (no P..AD coverage) ---
$ws.Range("F2").Value = 5
$ws.Range("P2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AD2").ClearContents()

# --- Refresh the view: zoom to 90%, reset scroll position, select AD2 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AD2").Select()
